$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap "Illes Balears" / "Illes Balears*" (rows 26 and 27)
$ws.Range("A26").Value = "Illes Balears*"
$ws.Range("A27").Value = "Illes Balears"

# Swap "Melilla" / "Huelva" (rows 52 and 53)
$ws.Range("A52").Value = "Huelva"
$ws.Range("A53").Value = "Melilla"

# Update the "last updated" timestamp text
$ws.Range("A1").Value = "Datos actualizados a 20 de Marzo de 2020 a las 00:15"
